$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Widen column B (bestFit column) to fit the new, longer "Question" text.
#    27.2421875 is the target character-width; the engine quantizes this to
#    the nearest whole pixel when applied through ColumnWidth.
$ws.Columns.Item(2).ColumnWidth = 27.2421875

# 2. Row 4 (A4) switches from the old date style (s=3) to the same date style
#    used by rows 2-3 (s=1). Copy the formatting from A2 and then restore A4's
#    original value (46066 = 2026-02-13), since Copy also copies the value.
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 46066

# 3. Add the new row 5 entry: "Min Opns. to convert w1 to w2" / Edit Distance.
# 3a. Date cell A5, formatted the same as the other date cells.
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 46067

# 3b. Question cell B5.
$ws.Range("B5").Value = "Min Opns. to convert w1 to w2"

# 3c. URL cell C5 with a real hyperlink, formatted like the other URL cells.
$ws.Hyperlinks.Add($ws.Range("C5"), "https://leetcode.com/problems/edit-distance/")
$ws.Range("C4").Copy($ws.Range("C5"))
$ws.Range("C5").Value = "https://leetcode.com/problems/edit-distance/"

# Adding the hyperlink registers a new built-in "Hyperlink" cell style that no
# cell actually ends up using (C5 keeps the pre-existing URL-column style
# copied from C4 above) - remove it so it doesn't linger unused.
$wb.Styles.Item("Hyperlink").Delete()
